$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.26437026758141258,
    -0.0099999997318676037,
    -0.0089999997303813473,
    0.061996648275520982,
    -0.0059999997375115299,
    -0.0059999997306334762,
    -0.019999999675153646,
    -0.019999999672027258,
    -0.0059999997229995827,
    -0.0059999997200748112,
    -0.0044999997257235691,
    -0.0059999997192288212,
    -0.0059999997167921038,
    0.010087515200567942,
    -0.0059999997150699258,
    -0.0059999997140209871,
    -0.0059999997126656268,
    -0.0089999997006815491,
    -0.058680798078850405,
    -0.0089999997324063941,
    -0.0089999997320111547,
    -0.0089999997317127267,
    -0.0089999997277807608,
    -0.041999999595007331,
    -0.041999999592715831,
    -0.005999999730587291,
    -0.0059999997302320196,
    -0.0059999997292825569,
    -0.011999999705606612,
    -0.019999999674731761,
    -0.014999999695717747,
    -0.02099999967322308,
    0.061383653262995885
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
